$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '21.109.51'
$ws.Range("E2").Value = '  +3.51%  '

$ws.Range("D3").Value = '1.535.78'
$ws.Range("E3").Value = '  +5.13%  '

$ws.Range("D4").Value = "'1.012"
$ws.Range("E4").Value = '  +0.21%  '

$ws.Range("D5").Value = "'0.9648"
$ws.Range("E5").Value = '  +1.54%  '

$ws.Range("D6").Value = "'281.92"
$ws.Range("E6").Value = '  +2.69%  '

$ws.Range("D7").Value = "'0.3623"
$ws.Range("E7").Value = '  -0.63%  '

$ws.Range("D8").Value = "'0.3178"
$ws.Range("E8").Value = '  +3.78%  '

$ws.Range("D9").Value = "'40.64"
$ws.Range("E9").Value = '  +2.20%  '

$ws.Range("D10").Value = "'1.094"
$ws.Range("E10").Value = '  +5.84%  '

$ws.Range("D11").Value = "'0.06801"
$ws.Range("E11").Value = '  +3.49%  '

$ws.Range("D12").Value = "'1.006"
$ws.Range("E12").Value = '  +0.40%  '

$ws.Range("D13").Value = "'5.663"
$ws.Range("E13").Value = '  +4.60%  '

$ws.Range("D14").Value = "'18.69"
$ws.Range("E14").Value = '  +4.25%  '

$ws.Range("D15").Value = "'6.344"
$ws.Range("E15").Value = '  +3.33%  '

$ws.Range("D16").Value = "'0.00001043"
$ws.Range("E16").Value = '  +2.12%  '

$ws.Range("D17").Value = "'0.9639"
$ws.Range("E17").Value = '  -0.68%  '

$ws.Range("D18").Value = '1.526.85'
$ws.Range("E18").Value = '  +4.45%  '

$ws.Range("D19").Value = "'0.06087"
$ws.Range("E19").Value = '  +4.82%  '

$ws.Range("D20").Value = "'72.18"
$ws.Range("E20").Value = '  +4.52%  '

$ws.Range("D21").Value = "'5.695"
$ws.Range("E21").Value = '  +4.88%  '

$ws.Range("D22").Value = "'14.98"
$ws.Range("E22").Value = '  +4.06%  '

$ws.Range("E23").Value = '  +4.24%  '

$ws.Range("D24").Value = "'2.327"
$ws.Range("E24").Value = '  +3.76%  '

$ws.Range("D25").Value = '21.157.94'
$ws.Range("E25").Value = '  +3.64%  '

$ws.Range("D26").Value = "'148.04"
$ws.Range("E26").Value = '  +4.52%  '

$ws.Range("D27").Value = "'2.213"
$ws.Range("E27").Value = '  +6.66%  '

$ws.Range("D28").Value = "'17.62"
$ws.Range("E28").Value = '  +3.04%  '

$ws.Range("D29").Value = '1.693.49'
$ws.Range("E29").Value = '  +4.81%  '

$ws.Range("D30").Value = "'118.31"
$ws.Range("E30").Value = '  +5.22%  '

$ws.Range("D31").Value = "'4.024"
$ws.Range("E31").Value = '  +4.33%  '

$ws.Range("D32").Value = "'0.8507"
$ws.Range("E32").Value = '  +7.81%  '

$ws.Range("E33").Value = '  +6.01%  '

$ws.Range("D34").Value = "'0.07985"
$ws.Range("E34").Value = '  +1.28%  '

$ws.Range("D35").Value = "'1.507"
$ws.Range("E35").Value = '  -1.13%  '

$ws.Range("D36").Value = "'4.958"
$ws.Range("E36").Value = '  +6.14%  '

$ws.Range("D37").Value = "'1.198"
$ws.Range("E37").Value = '  +4.39%  '

$ws.Range("D38").Value = "'0.05852"
$ws.Range("E38").Value = '  +2.68%  '

$ws.Range("E39").Value = '  +3.85%  '

$ws.Range("D40").Value = "'10.66"
$ws.Range("E40").Value = '  +3.40%  '

$ws.Range("D41").Value = "'7.715"
$ws.Range("E41").Value = '  +3.51%  '

$ws.Range("B42").Value = 'Frax'
$ws.Range("C42").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D42").Value = "'0.9643"
$ws.Range("E42").Value = '  +0.62%  '

$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").Value = "'0.1910"
$ws.Range("E43").Value = '  +2.98%  '

$ws.Range("D44").Value = "'0.5432"
$ws.Range("E44").Value = '  +3.41%  '

$ws.Range("D45").Value = "'12.56"
$ws.Range("E45").Value = '  +5.67%  '

$ws.Range("D46").Value = "'3.581"
$ws.Range("E46").Value = '  +2.72%  '

$ws.Range("D47").Value = "'0.5439"
$ws.Range("E47").Value = '  +5.77%  '

$ws.Range("D48").Value = "'121.48"
$ws.Range("E48").Value = '  +3.84%  '

$ws.Range("D49").Value = "'1.867"
$ws.Range("E49").Value = '  +7.12%  '

$ws.Range("D50").Value = "'0.06569"
$ws.Range("E50").Value = '  +2.42%  '

$ws.Range("D51").Value = "'0.9920"
$ws.Range("E51").Value = '  -0.06%  '
